# refatoracao e add excelConsumer
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CadastrarNovoCliente_Ne")
$ws.Activate()

# D2: Abc123 -> Abc1234
$ws.Range("D2").Value = "Abc1234"

# E4: the cell keeps its existing mailto hyperlink, but its displayed text is
# changed to "teste.com" while the hyperlink's original display name
# ("teste1@teste.com") is preserved.
$ws.Range("E4").Hyperlinks.Item(1).TextToDisplay = "teste1@teste.com"
$ws.Range("E4").Value = "teste.com"

# Move the active selection to B7 (single cell) before save.
$ws.Range("B7").Select()

$wb.Save()
